$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue $ws 'D2' '292.02'
Set-TextValue $ws 'E2' '-6.13%'
Set-TextValue $ws 'D3' '40.59'
Set-TextValue $ws 'E3' '1.75%'
Set-TextValue $ws 'D4' '5.017'
Set-TextValue $ws 'E4' '-1.57%'
Set-TextValue $ws 'D5' '0.07337'
Set-TextValue $ws 'E5' '-3.10%'
Set-TextValue $ws 'D6' '4.296'
Set-TextValue $ws 'E6' '-0.12%'
Set-TextValue $ws 'D7' '1.543'
Set-TextValue $ws 'E7' '-7.06%'
Set-TextValue $ws 'D8' '0.9223'
Set-TextValue $ws 'E8' '-0.74%'
$ws.Range('B9').Value = 'BTSEToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue $ws 'D9' '2.400'
Set-TextValue $ws 'E9' '-0.95%'
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws 'D10' '0.1218'
Set-TextValue $ws 'E10' '0.43%'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws 'D11' '0.1741'
Set-TextValue $ws 'E11' '-4.06%'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws 'D12' '0.08657'
Set-TextValue $ws 'E12' '-3.82%'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws 'D13' '0.04283'
Set-TextValue $ws 'E13' '2.98%'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws 'D14' '0.1052'
Set-TextValue $ws 'E14' '-0.21%'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws 'D15' '0.001275'
Set-TextValue $ws 'E15' '-0.57%'
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws 'D16' '0.005771'
Set-TextValue $ws 'E16' '0.53%'
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws 'D17' '3.339'
Set-TextValue $ws 'E17' '-0.34%'
Set-TextValue $ws 'D18' '0.3288'
Set-TextValue $ws 'E18' '-2.00%'
Set-TextValue $ws 'D19' '7.658'
Set-TextValue $ws 'E19' '0.05%'
Set-TextValue $ws 'E20' '2.80%'
Set-TextValue $ws 'D21' '0.2749'
Set-TextValue $ws 'E21' '-2.27%'
Set-TextValue $ws 'D22' '0.03935'
Set-TextValue $ws 'E22' '-2.16%'
Set-TextValue $ws 'E23' '-0.68%'
Set-TextValue $ws 'D24' '0.003776'
Set-TextValue $ws 'E24' '-7.84%'
Set-TextValue $ws 'E25' '0.79%'
Set-TextValue $ws 'D26' '0.0003728'
Set-TextValue $ws 'E26' '-95.04%'
Set-TextValue $ws 'D38' '0.02292'
Set-TextValue $ws 'E38' '-5.42%'
Set-TextValue $ws 'D39' '0.04989'
Set-TextValue $ws 'E39' '-3.09%'
Set-TextValue $ws 'D40' '0.005792'
Set-TextValue $ws 'E40' '164.82%'
Set-TextValue $ws 'D41' '0.007650'
Set-TextValue $ws 'E41' '-1.32%'
Set-TextValue $ws 'E42' '-1.05%'
Set-TextValue $ws 'D43' '0.007354'
Set-TextValue $ws 'E43' '-4.00%'
Set-TextValue $ws 'D44' '0.007783'
Set-TextValue $ws 'E44' '-2.55%'
Set-TextValue $ws 'D45' '0.3181'
Set-TextValue $ws 'E45' '2.26%'
Set-TextValue $ws 'D46' '0.00006370'
Set-TextValue $ws 'E46' '-3.39%'
Set-TextValue $ws 'E47' '-0.03%'
Set-TextValue $ws 'D48' '0.02045'
Set-TextValue $ws 'E48' '-93.21%'
Set-TextValue $ws 'D49' '0.00002103'
Set-TextValue $ws 'E49' '-0.03%'
Set-TextValue $ws 'D50' '0.0002003'
Set-TextValue $ws 'E50' '-0.03%'
